# Update the dSF column (F) values to reflect the re-pulled / re-computed
# data as described in the commit message ("repull data, push all data,
# mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -1
$ws.Range("F5").Value = -3
$ws.Range("F7").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("F14").Value = 0
